$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$blockCH = New-Object "object[,]" 24,6
$blockCH[0,0] = 4.920258781565997
$blockCH[0,1] = 9.456804409284379
$blockCH[0,2] = 13.00677307600122
$blockCH[0,3] = 29.13633053754067
$blockCH[0,4] = 30.65898060472242
$blockCH[0,5] = 14.04511577966951
$blockCH[1,0] = 4.749570390413714
$blockCH[1,1] = 9.493174656124046
$blockCH[1,2] = 13.08385970225123
$blockCH[1,3] = 29.00324814388141
$blockCH[1,4] = 30.13621996917638
$blockCH[1,5] = 14.03520142760284
$blockCH[2,0] = 4.643142525647408
$blockCH[2,1] = 9.51674282059877
$blockCH[2,2] = 13.13350793028787
$blockCH[2,3] = 28.93048671829183
$blockCH[2,4] = 29.82228752659161
$blockCH[2,5] = 14.03241551912754
$blockCH[3,0] = 4.599441837565622
$blockCH[3,1] = 9.526658541742643
$blockCH[3,2] = 13.15432401779487
$blockCH[3,3] = 28.90310860666817
$blockCH[3,4] = 29.69632018270886
$blockCH[3,5] = 14.0321111455634
$blockCH[4,0] = 4.592167812978447
$blockCH[4,1] = 9.528323869328071
$blockCH[4,2] = 13.15781583066772
$blockCH[4,3] = 28.89870032562276
$blockCH[4,4] = 29.67552745943667
$blockCH[4,5] = 14.03211079857332
$blockCH[5,0] = 4.642554393276916
$blockCH[5,1] = 9.516875285513949
$blockCH[5,2] = 13.13378629666387
$blockCH[5,3] = 28.93010825921808
$blockCH[5,4] = 29.82058049017369
$blockCH[5,5] = 14.03240804947066
$blockCH[6,0] = 4.86178747149533
$blockCH[6,1] = 9.469088514590878
$blockCH[6,2] = 13.03287259745841
$blockCH[6,3] = 29.08859997776683
$blockCH[6,4] = 30.47738628421752
$blockCH[6,5] = 14.04101266791539
$blockCH[7,0] = 5.275705716908866
$blockCH[7,1] = 9.385166978503563
$blockCH[7,2] = 12.85329498648699
$blockCH[7,3] = 29.46930115143144
$blockCH[7,4] = 31.8127579178313
$blockCH[7,5] = 14.08402786312987
$blockCH[8,0] = 5.56635123786754
$blockCH[8,1] = 9.329442485732493
$blockCH[8,2] = 12.73242658261246
$blockCH[8,3] = 29.78986705667299
$blockCH[8,4] = 32.81132743154357
$blockCH[8,5] = 14.13145965070632
$blockCH[9,0] = 5.694993491010759
$blockCH[9,1] = 9.305372974545586
$blockCH[9,2] = 12.67982286666566
$blockCH[9,3] = 29.94413419318374
$blockCH[9,4] = 33.26702790115322
$blockCH[9,5] = 14.15643917136421
$blockCH[10,0] = 5.74314636469456
$blockCH[10,1] = 9.296442009745807
$blockCH[10,2] = 12.66024393052381
$blockCH[10,3] = 30.00372337919578
$blockCH[10,4] = 33.43960725927923
$blockCH[10,5] = 14.16638351886989
$blockCH[11,0] = 5.732801432684322
$blockCH[11,1] = 9.298357291810957
$blockCH[11,2] = 12.66444546112955
$blockCH[11,3] = 29.99083837355332
$blockCH[11,4] = 33.40244150756456
$blockCH[11,5] = 14.16422032103952
$blockCH[12,0] = 5.698966554170629
$blockCH[12,1] = 9.304634540609134
$blockCH[12,2] = 12.67820526994892
$blockCH[12,3] = 29.94901340168991
$blockCH[12,4] = 33.28122687601986
$blockCH[12,5] = 14.15724760163679
$blockCH[13,0] = 5.678167328019181
$blockCH[13,1] = 9.308503439913233
$blockCH[13,2] = 12.68667791410779
$blockCH[13,3] = 29.92354567798551
$blockCH[13,4] = 33.20697583709445
$blockCH[13,5] = 14.15303965675957
$blockCH[14,0] = 5.557868082716372
$blockCH[14,1] = 9.331041203862748
$blockCH[14,2] = 12.73591213946646
$blockCH[14,3] = 29.7799516333115
$blockCH[14,4] = 32.78155882533033
$blockCH[14,5] = 14.12989530850392
$blockCH[15,0] = 5.483117637272674
$blockCH[15,1] = 9.345194899639141
$blockCH[15,2] = 12.76672438144677
$blockCH[15,3] = 29.69399393347192
$blockCH[15,4] = 32.52081059557507
$blockCH[15,5] = 14.11656572563363
$blockCH[16,0] = 5.439789224165495
$blockCH[16,1] = 9.353456219063801
$blockCH[16,2] = 12.78467085354397
$blockCH[16,3] = 29.64535045599834
$blockCH[16,4] = 32.37098235675436
$blockCH[16,5] = 14.10921945368963
$blockCH[17,0] = 5.425063146152888
$blockCH[17,1] = 9.356274066256717
$blockCH[17,2] = 12.79078574687187
$blockCH[17,3] = 29.62901875593693
$blockCH[17,4] = 32.32028435500107
$blockCH[17,5] = 14.10678730353628
$blockCH[18,0] = 5.49110987067779
$blockCH[18,1] = 9.343675747461075
$blockCH[18,2] = 12.76342118299855
$blockCH[18,3] = 29.70306207529817
$blockCH[18,4] = 32.54855386643516
$blockCH[18,5] = 14.1179515346278
$blockCH[19,0] = 5.708920256742102
$blockCH[19,1] = 9.302785779780027
$blockCH[19,2] = 12.67415443741015
$blockCH[19,3] = 29.96126696200672
$blockCH[19,4] = 33.31683163048842
$blockCH[19,5] = 14.15928252761404
$blockCH[20,0] = 5.84798300960249
$blockCH[20,1] = 9.277132018664238
$blockCH[20,2] = 12.61780019092449
$blockCH[20,3] = 30.13682630409719
$blockCH[20,4] = 33.81895171090304
$blockCH[20,5] = 14.18912017596336
$blockCH[21,0] = 5.774077729425549
$blockCH[21,1] = 9.290726120818741
$blockCH[21,2] = 12.64769614557191
$blockCH[21,3] = 30.04251865034409
$blockCH[21,4] = 33.55102192996768
$blockCH[21,5] = 14.1729382545329
$blockCH[22,0] = 5.487497679073073
$blockCH[22,1] = 9.344362169343531
$blockCH[22,2] = 12.76491383589809
$blockCH[22,3] = 29.69895995181912
$blockCH[22,4] = 32.53601087152747
$blockCH[22,5] = 14.11732402231599
$blockCH[23,0] = 5.165850155137005
$blockCH[23,1] = 9.406825689841504
$blockCH[23,2] = 12.89992498689189
$blockCH[23,3] = 29.35899542661668
$blockCH[23,4] = 31.4476177274501
$blockCH[23,5] = 14.06960123124142
$ws.Range("C2:H25").Value = $blockCH

$blockJ = New-Object "object[,]" 24,1
$blockJ[0,0] = 9.437723190109452
$blockJ[1,0] = 9.482599112501354
$blockJ[2,0] = 9.511466432683115
$blockJ[3,0] = 9.523561310495268
$blockJ[4,0] = 9.525589688387241
$blockJ[5,0] = 9.51162820611013
$blockJ[6,0] = 9.452924411932964
$blockJ[7,0] = 9.348182252552908
$blockJ[8,0] = 9.277489692695589
$blockJ[9,0] = 9.246675775488161
$blockJ[10,0] = 9.235199667813374
$blockJ[11,0] = 9.237662707450751
$blockJ[12,0] = 9.245727777162793
$blockJ[13,0] = 9.25069290238685
$blockJ[14,0] = 9.279530433635543
$blockJ[15,0] = 9.297565050561989
$blockJ[16,0] = 9.308064686536371
$blockJ[17,0] = 9.311641453728463
$blockJ[18,0] = 9.295632138427617
$blockJ[19,0] = 9.243353654531145
$blockJ[20,0] = 9.210308008057137
$blockJ[21,0] = 9.227842776048046
$blockJ[22,0] = 9.296505598936786
$blockJ[23,0] = 9.348182252552908
$ws.Range("J2:J25").Value = $blockJ

$blockMO = New-Object "object[,]" 24,3
$blockMO[0,0] = 22.35015264789025
$blockMO[0,1] = 17.43041460879761
$blockMO[0,2] = 21.78685146966273
$blockMO[1,0] = 21.66593927181829
$blockMO[1,1] = 17.14105371264457
$blockMO[1,2] = 21.67356376577764
$blockMO[2,0] = 21.23465837880838
$blockMO[2,1] = 16.9628288749524
$blockMO[2,2] = 21.61021082340978
$blockMO[3,0] = 21.05633064869456
$blockMO[3,1] = 16.89014992225147
$blockMO[3,2] = 21.58597489955432
$blockMO[4,0] = 21.02657095312731
$blockMO[4,1] = 16.87808119386683
$blockMO[4,2] = 21.58204656828483
$blockMO[5,0] = 21.2322635036619
$blockMO[5,1] = 16.96184878854882
$blockMO[5,2] = 21.60987754337788
$blockMO[6,0] = 22.11668759427191
$blockMO[6,1] = 17.33081122424798
$blockMO[6,2] = 21.746514767159
$blockMO[7,0] = 23.75318926564491
$blockMO[7,1] = 18.04612151014032
$blockMO[7,2] = 22.06268758165521
$blockMO[8,0] = 24.88486897604502
$blockMO[8,1] = 18.56171446959201
$blockMO[8,2] = 22.32285378145962
$blockMO[9,0] = 25.38239477384256
$blockMO[9,1] = 18.79312710310437
$blockMO[9,2] = 22.44689102276475
$blockMO[10,0] = 25.56816910677947
$blockMO[10,1] = 18.88022915096149
$blockMO[10,2] = 22.49464417921837
$blockMO[11,0] = 25.52827826385139
$blockMO[11,1] = 18.86149489446174
$blockMO[11,2] = 22.48432543161342
$blockMO[12,0] = 25.39773196912811
$blockMO[12,1] = 18.80030407549908
$blockMO[12,2] = 22.45080421720334
$blockMO[13,0] = 25.31742239505121
$blockMO[13,1] = 18.76275179671525
$blockMO[13,2] = 22.4303724084171
$blockMO[14,0] = 24.85199417630853
$blockMO[14,1] = 18.54652127515741
$blockMO[14,2] = 22.31485914856921
$blockMO[15,0] = 24.56193617281236
$blockMO[15,1] = 18.41301138737857
$blockMO[15,2] = 22.24542774120273
$blockMO[16,0] = 24.39348340110551
$blockMO[16,1] = 18.33592956422917
$blockMO[16,2] = 22.20603005163185
$blockMO[17,0] = 24.336174725997
$blockMO[17,1] = 18.3097835466969
$blockMO[17,2] = 22.19278400373203
$blockMO[18,0] = 24.592981996845
$blockMO[18,1] = 18.42725438507583
$blockMO[18,2] = 22.25276345837043
$blockMO[19,0] = 25.43614896094802
$blockMO[19,1] = 18.81829226107818
$blockMO[19,2] = 22.46062924664157
$blockMO[20,0] = 25.97183461631777
$blockMO[20,1] = 19.07074180755329
$blockMO[20,2] = 22.60102792116418
$blockMO[21,0] = 25.68737887540235
$blockMO[21,1] = 18.93631472751926
$blockMO[21,2] = 22.52569040917387
$blockMO[22,0] = 24.5789514469152
$blockMO[22,1] = 18.42081613276582
$blockMO[22,2] = 22.24944536018794
$blockMO[23,0] = 23.75318926564491
$blockMO[23,1] = 17.85401202175196
$blockMO[23,2] = 22.06268758165521
$ws.Range("M2:O25").Value = $blockMO

Write-Host "Updated loading_percent values for rows 2-25 (case with 380 kV)"
